$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 updates
$ws.Range("F2").Value = 1.31
$ws.Range("G2").Value = 1.35
$ws.Range("H2").Value = 12
$ws.Range("I2").Value = 14
$ws.Range("J2").Value = 5.3
$ws.Range("K2").Value = 6.2
$ws.Range("L2").Value = 1.28
$ws.Range("N2").Value = 4.1
$ws.Range("O2").Value = 1.28
$ws.Range("P2").Value = 2.08
$ws.Range("Q2").Value = 1.79
$ws.Range("R2").Value = 1.42
$ws.Range("S2").Value = 3.05
$ws.Range("T2").Value = 2.24
$ws.Range("U2").Value = 1.7
$ws.Range("V2").Value = 1.07
$ws.Range("W2").Value = 3.85
$ws.Range("X2").Value = 22
$ws.Range("Z2").Value = 140
$ws.Range("AE2").Value = 280
$ws.Range("AH2").Value = 1000
$ws.Range("AI2").Value = 210
$ws.Range("AL2").Value = 1000
$ws.Range("AM2").Value = 270
$ws.Range("AN2").Value = 6.2

# Row 4 updates
$ws.Range("I4").Value = 14
$ws.Range("J4").Value = 2.94

$wb.Save()
